$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values like "0.90%" or "0.06420" must be written as literal text
# (matching the source inlineStr cells), not auto-converted by Excel
# into numbers/percentages. Setting NumberFormat to Text ("@") first
# forces the literal string to be stored verbatim.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.90%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.10%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.888"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.02%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06420"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.26%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.970"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.14%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.184"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.46%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8839"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.86%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1536"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.70%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05146"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.27%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07412"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.03%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02887"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.66%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08975"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.65%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001564"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.87%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006366"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.34%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006154"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "5.52%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.480"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.84%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.314"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.09%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.08%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.44%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.902"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.43%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04423"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.56%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "8.72%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.001177"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.14%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-9.14%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.53%"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "15.72%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04142"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.87%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006788"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.33%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.61%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001901"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.34%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01187"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.99%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005312"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.07%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.685"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "13.10%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01853"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-7.26%"
